# Update gh-pages to output generated at 456a3b4
#
# Refreshes the "want-to-go" head-count (column F) for each con/event row
# across all four sheets, and flips two rows whose lowest ticket price
# (column G) sold out since the last crawl (numeric price -> "已售罄").

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value  = 902    # OVO动漫嘉年华2th: 901 -> 902
$ws.Cells.Item(6, 6).Value  = 351    # 原神&崩铁同人only: 352 -> 351
$ws.Cells.Item(7, 6).Value  = 389    # 虚拟主播Virtual Only: 390 -> 389
$ws.Cells.Item(7, 7).Value  = "已售罄"  # ...and its lowest price (80) sold out
$ws.Cells.Item(8, 6).Value  = 737    # 迷宫饭only: 738 -> 737
$ws.Cells.Item(9, 6).Value  = 1124   # 南国书香节璃樱动漫嘉年华: 1113 -> 1124
$ws.Cells.Item(10, 6).Value = 12446  # 第九届初物语动漫展: 12409 -> 12446
$ws.Cells.Item(11, 6).Value = 683    # COC星火次元云漫创作交流展: 680 -> 683
$ws.Cells.Item(13, 6).Value = 314    # 第九届初物语动漫展内场—羊仔: 313 -> 314
$ws.Cells.Item(16, 6).Value = 309    # 凹凸世界ONLY: 306 -> 309
$ws.Cells.Item(17, 6).Value = 1831   # 原神×崩坏×绝区零同人only: 1829 -> 1831
$ws.Cells.Item(20, 6).Value = 505    # LoveLiveOnly: 504 -> 505
$ws.Cells.Item(22, 6).Value = 118    # 第九届初物语动漫展内场—赵路: 117 -> 118
$ws.Cells.Item(24, 6).Value = 218    # 原神X崩坏X星铁旅行盛宴同人only: 216 -> 218
$ws.Cells.Item(26, 6).Value = 105    # 第七届AP动漫嘉年华: 102 -> 105
$ws.Cells.Item(27, 6).Value = 112    # 星光次元动漫嘉年华: 109 -> 112
$ws.Cells.Item(29, 6).Value = 198    # 樱漫潮玩动漫嘉年华: 194 -> 198
$ws.Cells.Item(30, 6).Value = 229    # wio流金序曲乙女同人展: 227 -> 229
$ws.Cells.Item(32, 6).Value = 53     # wio jumponly4.0万圣狂欢节: 52 -> 53

# ---- Sheet: 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value  = 34     # 33 -> 34
$ws.Cells.Item(3, 6).Value  = 165    # 166 -> 165
$ws.Cells.Item(6, 6).Value  = 267    # 266 -> 267
$ws.Cells.Item(7, 6).Value  = 4446   # 4443 -> 4446
$ws.Cells.Item(8, 6).Value  = 118    # 115 -> 118
$ws.Cells.Item(12, 6).Value = 339    # 337 -> 339

# ---- Sheet: 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 849     # 847 -> 849

# ---- Sheet: 全部类型 (All types, union of the sheets above) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 849    # (本地生活 row) 847 -> 849
$ws.Cells.Item(4, 6).Value  = 165    # (演出 row) 166 -> 165
$ws.Cells.Item(5, 6).Value  = 902    # (展览 row) 901 -> 902
$ws.Cells.Item(8, 6).Value  = 389    # (展览 row) 390 -> 389
$ws.Cells.Item(8, 7).Value  = "已售罄"  # ...and its lowest price (80) sold out
$ws.Cells.Item(11, 6).Value = 1124   # 1113 -> 1124
$ws.Cells.Item(12, 6).Value = 12446  # 12409 -> 12446
$ws.Cells.Item(13, 6).Value = 267    # 266 -> 267
$ws.Cells.Item(14, 6).Value = 683    # 680 -> 683
$ws.Cells.Item(16, 6).Value = 314    # 313 -> 314
$ws.Cells.Item(18, 6).Value = 309    # 306 -> 309
$ws.Cells.Item(19, 6).Value = 1831   # 1829 -> 1831
$ws.Cells.Item(22, 6).Value = 505    # 504 -> 505
$ws.Cells.Item(23, 6).Value = 4446   # 4443 -> 4446
$ws.Cells.Item(25, 6).Value = 118    # 115 -> 118
$ws.Cells.Item(26, 6).Value = 118    # 115 -> 118
$ws.Cells.Item(28, 6).Value = 118    # 117 -> 118
$ws.Cells.Item(31, 6).Value = 339    # 337 -> 339
$ws.Cells.Item(34, 6).Value = 218    # 216 -> 218
$ws.Cells.Item(36, 6).Value = 105    # 102 -> 105
$ws.Cells.Item(37, 6).Value = 112    # 109 -> 112
$ws.Cells.Item(40, 6).Value = 198    # 194 -> 198
$ws.Cells.Item(43, 6).Value = 229    # 227 -> 229
$ws.Cells.Item(46, 6).Value = 54     # 52 -> 54
